$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the notes text in D9 - append the new sentence about the Adventurer Agent class
$ws.Range("D9").Value = "Added Wumpus and stenches in addition to Adventurer. Renders well, simplify the iconography so I can just use core drawing functions. Solid ovals for things, empty ovals for their signals. Must update documention. Got the rest of the images drawn for the pits, breezes, gold and glitter. Need to turn of the random seed as we get no variation in runs. Created the Adventurer Agent class that initially moves randomly."

# Update hours spent for the day (C9) from 4.5 to 5.5
$ws.Range("C9").Value = 5.5

# Update the selected cell to D10
$ws.Range("D10").Select()

# Update row 9 height to account for the longer wrapped text
$ws.Rows.Item(9).RowHeight = 71.25
